# Atualizado por script em 26-11-2023 20:30
#
# 1) A handful of existing match rows had their home/away sides (and all the
#    associated odds/timestamps/url columns, F:V) swapped with another row
#    on the sheet. Columns A:E (index, pais, torneio, temporada, data_partida)
#    are untouched for every row.
# 2) Five new match rows (sheet rows 52-56) are appended at the bottom with
#    the next sequential Indice values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($Sheet, $RowA, $RowB, $FirstCol, $LastCol)
    for ($c = $FirstCol; $c -le $LastCol; $c++) {
        $valA = $Sheet.Cells.Item($RowA, $c).Value2
        $valB = $Sheet.Cells.Item($RowB, $c).Value2
        $Sheet.Cells.Item($RowA, $c).Value2 = $valB
        $Sheet.Cells.Item($RowB, $c).Value2 = $valA
    }
}

# Swap F:V (columns 6-22) between the paired rows.
Swap-RowRange $ws 14 15 6 22
Swap-RowRange $ws 25 26 6 22
Swap-RowRange $ws 38 40 6 22
Swap-RowRange $ws 39 41 6 22
Swap-RowRange $ws 43 45 6 22

function Set-MatchRow {
    param(
        $Sheet, $Row, $Indice, $DataPartida,
        $Home, $HomeGols, $Away, $AwayGols,
        $HomeOpen, $HomeOpenDt, $HomeClose, $HomeCloseDt,
        $DrawOpen, $DrawOpenDt, $DrawClose, $DrawCloseDt,
        $AwayOpen, $AwayOpenDt, $AwayClose, $AwayCloseDt,
        $Url
    )

    # Pull the A/E column formatting (bold+border+centered / date number
    # format) straight from the template row above instead of re-declaring
    # fonts/borders, so no new style entries get created.
    $Sheet.Range("A51").Copy()
    $Sheet.Cells.Item($Row, 1).PasteSpecial(-4122)
    $Sheet.Range("E51").Copy()
    $Sheet.Cells.Item($Row, 5).PasteSpecial(-4122)

    $Sheet.Cells.Item($Row, 1).Value2 = $Indice

    $Sheet.Cells.Item($Row, 2).Value2 = "cambodia"
    $Sheet.Cells.Item($Row, 3).Value2 = "cpl"
    $Sheet.Cells.Item($Row, 4).Value2 = "2023-2024"

    $Sheet.Cells.Item($Row, 5).Value2 = $DataPartida

    $Sheet.Cells.Item($Row, 6).Value2 = $Home
    $Sheet.Cells.Item($Row, 7).Value2 = $HomeGols
    $Sheet.Cells.Item($Row, 8).Value2 = $Away
    $Sheet.Cells.Item($Row, 9).Value2 = $AwayGols

    $Sheet.Cells.Item($Row, 10).Value2 = $HomeOpen
    $Sheet.Cells.Item($Row, 11).Value2 = $HomeOpenDt
    $Sheet.Cells.Item($Row, 12).Value2 = $HomeClose
    $Sheet.Cells.Item($Row, 13).Value2 = $HomeCloseDt

    $Sheet.Cells.Item($Row, 14).Value2 = $DrawOpen
    $Sheet.Cells.Item($Row, 15).Value2 = $DrawOpenDt
    $Sheet.Cells.Item($Row, 16).Value2 = $DrawClose
    $Sheet.Cells.Item($Row, 17).Value2 = $DrawCloseDt

    $Sheet.Cells.Item($Row, 18).Value2 = $AwayOpen
    $Sheet.Cells.Item($Row, 19).Value2 = $AwayOpenDt
    $Sheet.Cells.Item($Row, 20).Value2 = $AwayClose
    $Sheet.Cells.Item($Row, 21).Value2 = $AwayCloseDt

    $Sheet.Cells.Item($Row, 22).Value2 = $Url
}

Set-MatchRow $ws 52 51 45255.39583333334 `
    "Prey Veng" 0 "Dangkor" 1 `
    1.71 "24/11/2023 22:43" 1.79 "25/11/2023 01:42" `
    3.97 "24/11/2023 22:43" 4.02 "25/11/2023 09:04" `
    3.63 "24/11/2023 22:43" 3.33 "25/11/2023 09:04" `
    "https://www.betexplorer.com/football/cambodia/cpl/prey-veng-dangkor-senchey/nqo8zb7D/"

Set-MatchRow $ws 53 52 45255.5 `
    "Phnom Penh Crown" 0 "Svay Rieng" 0 `
    1.93 "25/11/2023 00:12" 2.11 "25/11/2023 11:58" `
    3.55 "25/11/2023 00:12" 3.53 "25/11/2023 11:58" `
    3.22 "25/11/2023 00:12" 2.88 "25/11/2023 11:58" `
    "https://www.betexplorer.com/football/cambodia/cpl/phnom-penh-crown-svay-rieng/dESbcLEQ/"

Set-MatchRow $ws 54 53 45255.5 `
    "Visakha" 1 "NagaWorld" 1 `
    1.37 "25/11/2023 00:12" 1.37 "25/11/2023 00:16" `
    4.73 "25/11/2023 00:12" 4.82 "25/11/2023 11:56" `
    5.85 "25/11/2023 00:12" 5.98 "25/11/2023 11:56" `
    "https://www.betexplorer.com/football/cambodia/cpl/visakha-nagaworld/z9s4yvh7/"

Set-MatchRow $ws 55 54 45256.39583333334 `
    "Angkor Tiger" 3 "Boeung Ket" 3 `
    4.26 "25/11/2023 22:42" 5.06 "26/11/2023 09:20" `
    4.18 "25/11/2023 22:42" 4.49 "26/11/2023 09:20" `
    1.54 "25/11/2023 22:42" 1.46 "26/11/2023 09:20" `
    "https://www.betexplorer.com/football/cambodia/cpl/angkor-tiger-boeung-ket/vwYgb1aK/"

Set-MatchRow $ws 56 55 45256.5 `
    "Tiffy Army" 4 "Kirivong Sok Sen Chey" 2 `
    1.57 "26/11/2023 00:42" 1.76 "26/11/2023 11:32" `
    3.76 "26/11/2023 00:42" 3.56 "26/11/2023 11:32" `
    4.56 "26/11/2023 00:42" 3.88 "26/11/2023 11:32" `
    "https://www.betexplorer.com/football/cambodia/cpl/tiffy-army-kirivong-sok-sen-chey/0nZkaspE/"
